# link updated on 20190308
# Adds a new "Tipo indirizzo" lookup row (11 -> "Rettifica post accertamenti")
# to the bottom of the table, then refreshes column B's width to fit the
# text and leaves the selection on the first empty row below the table,
# mirroring the manual data-entry workflow that produced the upstream commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New last row of the table: ID 11, DESCRIZIONE "Rettifica post accertamenti"
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = "Rettifica post accertamenti"

# Resize column B so the new (longer) description fits ("best fit" width)
$ws.Columns.Item(2).ColumnWidth = 29.5

# Move the active selection to the row right after the new data, like a
# user would after typing the new entry and pressing Enter/Down.
$ws.Range("A13").Select()
